$d = $word.ActiveDocument

# 1) Replace the two recurring tuple lines (5 occurrences each, throughout the document).
$r1 = $d.Content
[void]$r1.Find.Execute(
    "(Entity, Statement, Attribute, Value);", $true, $false, $false, $false, $false,
    $true, 1, $false, "(Entity, Statement, Occurrence, Attribute);", 2)

$r2 = $d.Content
[void]$r2.Find.Execute(
    "(Role, Entity, Statement, Attribute);", $true, $false, $false, $false, $false,
    $true, 1, $false, "(Role, Entity, Statement, Occurrence);", 2)

# 2) Append a new block of paragraphs right after the "Dataflow contexts from Message
#    levels application." paragraph (in each place Find lands on it -- there is exactly
#    one occurrence in this document).
$r3 = $d.Content
[void]$r3.Find.Execute("Dataflow contexts from Message levels application.")
$r3.Collapse(0)

$newParas = @(
    $null,
    "Augmentation: For example, a template Statement (Statement used as transform specification) from, for example, the Interaction Model, may state matching pattetns such as:",
    $null,
    "(ContextClass : Subject, Context, Occurrence, Attribute);",
    $null,
    "and, when applied to an input Message:",
    $null,
    "(Statement, Subject, Predicate, Value);",
    $null,
    "reacts emitting the following Statement, transforming input context Message according template rules (input Subject -> output Attribute):",
    $null,
    "(TransformClass : Entity, Statement, Subject, Predicate);",
    $null,
    "which is materialized in the corresponding Model and is itself again a Message routed for further processing. TransformClass is an instance / subclass of super / meta class ContextClass (model layers transform rules)."
)

foreach ($t in $newParas) {
    [void]$r3.InsertParagraphAfter()
    [void]$r3.Move(1, 1)
    if ($t -ne $null) {
        [void]$r3.InsertAfter($t)
    }
}
